# Add 2022-Q3 data
# -----------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q3" right after "总计" by
#    copying the existing "2022-Q2" sheet (same column layout/styles)
#    and overwriting its data with the Q3 numbers.
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    at the top of the data and push the existing rows down by one.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet from a copy of "2022-Q2" ----
$q2.Copy($null, $summary)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# "2022-Q2" only had 3 data rows (rows 2-4); "2022-Q3" needs 5 (rows
# 2-6). Stretch the A-column style (bold/border/center, like A2:A4)
# down to A5:A6 before writing the new values.
$q3.Range("A4").Copy()
$q3.Range("A5:A6").PasteSpecial(-4122)

# Fund-code / numeric-looking text columns must stay TEXT (leading
# zeros, trailing zeros like "0.20" must be preserved) -> prefix the
# literal apostrophe so Excel stores them as strings, not numbers.
$q3.Range("A2").Value2 = 0
$q3.Range("B2").Value2 = "'004497"
$q3.Range("C2").Value2 = "'前海开源多元策略灵活配置混合C"
$q3.Range("D2").Value2 = "'1.68"
$q3.Range("E2").Value2 = "'93.04"
$q3.Range("F2").Value2 = "'8.85"
$q3.Range("G2").Value2 = "'0.1487"
$q3.Range("H2").Value2 = 1

$q3.Range("A3").Value2 = 1
$q3.Range("B3").Value2 = "'004496"
$q3.Range("C3").Value2 = "'前海开源多元策略灵活配置混合A"
$q3.Range("D3").Value2 = "'0.91"
$q3.Range("E3").Value2 = "'93.04"
$q3.Range("F3").Value2 = "'8.85"
$q3.Range("G3").Value2 = "'0.0805"
$q3.Range("H3").Value2 = 1

$q3.Range("A4").Value2 = 2
$q3.Range("B4").Value2 = "'003993"
$q3.Range("C4").Value2 = "'前海开源沪港深核心驱动灵活配置混合"
$q3.Range("D4").Value2 = "'0.53"
$q3.Range("E4").Value2 = "'82.41"
$q3.Range("F4").Value2 = "'8.42"
$q3.Range("G4").Value2 = "'0.0446"
$q3.Range("H4").Value2 = 2

$q3.Range("A5").Value2 = 3
$q3.Range("B5").Value2 = "'161124"
$q3.Range("C5").Value2 = "'易方达香港恒生综合小型股指数（QDII-LOF）A"
$q3.Range("D5").Value2 = "'0.20"
$q3.Range("E5").Value2 = "'91.61"
$q3.Range("F5").Value2 = "'1.47"
$q3.Range("G5").Value2 = "'0.0029"
$q3.Range("H5").Value2 = 5

$q3.Range("A6").Value2 = 4
$q3.Range("B6").Value2 = "'006263"
$q3.Range("C6").Value2 = "'易方达香港恒生综合小型股指数（QDII-LOF）C"
$q3.Range("D6").Value2 = "'0.05"
$q3.Range("E6").Value2 = "'91.61"
$q3.Range("F6").Value2 = "'1.47"
$q3.Range("G6").Value2 = "'0.0007"
$q3.Range("H6").Value2 = 5

# --- 2. Update the "总计" sheet: shift rows down, insert new top row
$summary.Range("A8").Value2 = 6
$summary.Range("B8").Value2 = $summary.Range("B7").Value2
$summary.Range("C8").Value2 = $summary.Range("C7").Value2
$summary.Range("D8").Value2 = $summary.Range("D7").Value2

$summary.Range("B7").Value2 = $summary.Range("B6").Value2
$summary.Range("C7").Value2 = $summary.Range("C6").Value2
$summary.Range("D7").Value2 = $summary.Range("D6").Value2

$summary.Range("B6").Value2 = $summary.Range("B5").Value2
$summary.Range("C6").Value2 = $summary.Range("C5").Value2
$summary.Range("D6").Value2 = $summary.Range("D5").Value2

$summary.Range("B5").Value2 = $summary.Range("B4").Value2
$summary.Range("C5").Value2 = $summary.Range("C4").Value2
$summary.Range("D5").Value2 = $summary.Range("D4").Value2

$summary.Range("B4").Value2 = $summary.Range("B3").Value2
$summary.Range("C4").Value2 = $summary.Range("C3").Value2
$summary.Range("D4").Value2 = $summary.Range("D3").Value2

$summary.Range("B3").Value2 = $summary.Range("B2").Value2
$summary.Range("C3").Value2 = $summary.Range("C2").Value2
$summary.Range("D3").Value2 = $summary.Range("D2").Value2

$summary.Range("B2").Value2 = "2022-Q3"
$summary.Range("C2").Value2 = 5
$summary.Range("D2").Value2 = 0.28

# A8 needs the same bold/border/center style as the other A-column
# header cells (A2:A7 already carry it).
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)
$summary.Range("A8").Value2 = 6

# --- keep the original last sheet ("2020-Q4") as the selected tab --
$wb.Worksheets.Item($wb.Worksheets.Count).Select()
